$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 40: add commit text and hours
$ws.Range("C40").Value = "arena win/lose conditions added"
$ws.Range("G40").Value = 2

# Row 41: add commit text and hours
$ws.Range("C41").Value = "RESET update & bugs fixed"
$ws.Range("G41").Value = 2.5

# Extend the total formula to include the new rows
$ws.Range("G49").Formula = "=SUM(G4:G41)"

# Update the view: scroll position and active selection
$ws.Range("C42").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
